# Applies the changes described by the commit:
# "Added more data tables in csv_exports, updated README, added project-update.md"
#
# Concretely, within this workbook the edit:
#  - adds new source-file-name values into the "Sheet2" worksheet (years table)
#  - widens a few columns on that sheet to fit the new text
#  - switches the active/selected sheet from "Sheet1" to "Sheet2"
#  - freezes the first column on "Sheet2" and scrolls/selects near column H

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")   # years table (sheet1.xml, uid ...84B5)
$ws1 = $wb.Worksheets.Item("Sheet1")   # data journal (sheet2.xml, uid ...BB5)

# ---------------------------------------------------------------------------
# 1. New shared-string values written into the years table on "Sheet2"
# ---------------------------------------------------------------------------

$felonyFile      = "ca_felony_arrests_2015-2020"
$miscGenderFile  = "ca_misdemeanor_arrests_offense_by_gender_and_race_ethnic_group_2020"
$miscAdultFile   = "ca_misdeamnor_arrests_.by_offense_for_adult_and_juvenile_arrests"
$felonyAgeFile   = "felony_arrests_category_and_offense_by_age_group_of_arrestee_2020"
$miscAgeFile     = "misdemeanor_arrests_offense_by_age_group_of_arrestee_2020.csv"
$felonyGenderFile= "felony_arrests_category_and_offense_by_gender_and_race_ethnic_group_of_arrestee._2020.csv"

# Row 2 (year 2020) gets a value in every column, B..I
$ws2.Range("B2").Value = $felonyFile
$ws2.Range("C2").Value = $miscGenderFile
$ws2.Range("D2").Value = $felonyFile
$ws2.Range("E2").Value = $miscAdultFile
$ws2.Range("F2").Value = $felonyAgeFile
$ws2.Range("G2").Value = $miscAgeFile
$ws2.Range("H2").Value = $felonyGenderFile
$ws2.Range("I2").Value = $miscGenderFile

# Rows 3-7 (years 2019-2015) only get the felony-arrests filename in B and D
foreach ($r in 3..7) {
    $ws2.Range("B$r").Value = $felonyFile
    $ws2.Range("D$r").Value = $felonyFile
}

# ---------------------------------------------------------------------------
# 2. Column width changes on "Sheet2" (C, D, H widened to fit new text)
# ---------------------------------------------------------------------------
# NOTE: this runtime stores column width as `ColumnWidth + 5/6` in the saved
# file, and rounds ColumnWidth itself to the nearest 1/6 of a character, so
# we pick the ColumnWidth value that lands closest to the target file width.

$ws2.Columns.Item(3).ColumnWidth = 45.5               # -> width ~46.33
$ws2.Columns.Item(4).ColumnWidth = 28.6666666666667    # -> width = 29.5
$ws2.Columns.Item(8).ColumnWidth = 44.3333333333333    # -> width ~45.16

# ---------------------------------------------------------------------------
# 3. Switch the selected/active sheet from "Sheet1" to "Sheet2"
# ---------------------------------------------------------------------------

$ws2.Activate()

# ---------------------------------------------------------------------------
# 4. Freeze the first column on "Sheet2" and select near column H
# ---------------------------------------------------------------------------

$ws2.Range("H3").Select()
$excel.ActiveWindow.SplitColumn = 1
$excel.ActiveWindow.SplitRow = 0
$excel.ActiveWindow.FreezePanes = $true
